$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for first row (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-20 15:05:59"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-20 15:05:55"
$wsZhCn.Range("K2").Value = "2016-08-20 15:06:14"

# de-de sheet: Correspond Handback DateTime (K2) -- H2 shares the same text as
# Overview!G2 ("2016-08-20 15:05:17") and is updated together with it above.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-20 15:05:59"
$wsDeDe.Range("K2").Value = "2016-08-20 15:06:20"
